# Update the "K" (strikeout) column (column G) of the save_data sheet.
# The commit regenerates save_data to use K instead of Strike#, recalculating
# the K values for each logged appearance (row 2 through row 42, column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$gUpdates = @{
    2  = 3
    3  = 0
    4  = 2
    5  = 0
    6  = 5
    7  = 3
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 3
    22 = 2
    23 = 0
    24 = 0
    25 = 2
    26 = 0
    27 = 2
    28 = 0
    29 = 3
    30 = 0
    31 = 0
    32 = 1
    33 = 3
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    39 = 0
    40 = 2
    42 = 1
}

foreach ($row in $gUpdates.Keys) {
    $ws.Cells.Item($row, 7).Value = $gUpdates[$row]
}
